# Weekly price update for "Fruta, Terminal La Palmera de La Serena - Frutilla":
# a new week's worth of data (3 rows: Especial / Primera / Segunda) is inserted
# at the top of the price history table (row 313), pushing all existing rows
# down by 3 (old 313:402 -> new 316:405).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the current row 313, shifting the rest of the
# table (old rows 313:402) down to 316:405. Formatting (e.g. the date style
# on column D) is inherited from the row above, matching the rest of the
# table.
$ws.Rows("313:315").Insert()

# New rows of data (date serial 44463 = 2021-09-24), matching the existing
# record layout used throughout this sheet.
$newRows = @(
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44463, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Especial", 160, 25000, 26000, 25500, "`$/bandeja 7 kilos", "Provincia de Melipilla", 3643, 7),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44463, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera", 200, 22000, 23000, 22500, "`$/bandeja 7 kilos", "Provincia de Melipilla", 3214, 7),
    @(8, "Terminal La Palmera de La Serena", "Coquimbo", 44463, 4, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Segunda", 200, 19000, 20000, 19500, "`$/bandeja 7 kilos", "Provincia de Melipilla", 2786, 7)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $rowNum = 313 + $i
    $rowValues = $newRows[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowValues[$col - 1]
    }
}
